$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row19 = $ws.Range("A19:C19")
$row21 = $ws.Range("A21:C21")
$scratch = $ws.Range("A100:C100")

# --- Swap the (direct) cell formatting between row 19 and row 21 ---
# Stash row21's current formatting in a scratch range.
$row21.Copy()
$scratch.PasteSpecial(-4122)  # xlPasteFormats

# Apply row19's (old) formatting onto row21.
$row19.Copy()
$row21.PasteSpecial(-4122)   # xlPasteFormats

# Apply the stashed (old row21) formatting onto row19.
$scratch.Copy()
$row19.PasteSpecial(-4122)   # xlPasteFormats

# Clean up the scratch range so it doesn't leak into the sheet.
$scratch.Clear()
$excel.CutCopyMode = 0

# --- Swap the Y/N values held in column C for rows 19 and 21 ---
$c19 = $ws.Range("C19").Value2
$c21 = $ws.Range("C21").Value2
$ws.Range("C19").Value2 = $c21
$ws.Range("C21").Value2 = $c19

# --- Move the active selection from B24 to B21 ---
$ws.Range("B21").Select()
